# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for Mandarina - Clemenuless right before the
# existing row 244, shifting the subsequent rows (old 244-263) down to 247-266.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above row 244 (old row 244 and below shift down to 247+)
$ws.Range("A244:T246").EntireRow.Insert()

# New row 244
$ws.Range("A244").Value = 5
$ws.Range("B244").Value = "Macroferia Regional de Talca"
$ws.Range("C244").Value = "Maule"
$ws.Range("D244").Value = 44714
$ws.Range("E244").Value = 7
$ws.Range("F244").Value = "Fruta"
$ws.Range("G244").Value = 100102
$ws.Range("H244").Value = "Cítricos"
$ws.Range("I244").Value = 100102004
$ws.Range("J244").Value = "Mandarina"
$ws.Range("K244").Value = "Clemenuless"
$ws.Range("L244").Value = "Primera"
$ws.Range("M244").Value = 230
$ws.Range("N244").Value = 7000
$ws.Range("O244").Value = 7000
$ws.Range("P244").Value = 7000
$ws.Range("Q244").Value = "`$/caja 10 kilos"
$ws.Range("R244").Value = "Provincia de Limarí"
$ws.Range("S244").Value = 700
$ws.Range("T244").Value = 10

# New row 245
$ws.Range("A245").Value = 5
$ws.Range("B245").Value = "Macroferia Regional de Talca"
$ws.Range("C245").Value = "Maule"
$ws.Range("D245").Value = 44714
$ws.Range("E245").Value = 7
$ws.Range("F245").Value = "Fruta"
$ws.Range("G245").Value = 100102
$ws.Range("H245").Value = "Cítricos"
$ws.Range("I245").Value = 100102004
$ws.Range("J245").Value = "Mandarina"
$ws.Range("K245").Value = "Clemenuless"
$ws.Range("L245").Value = "Primera"
$ws.Range("M245").Value = 230
$ws.Range("N245").Value = 10000
$ws.Range("O245").Value = 10000
$ws.Range("P245").Value = 10000
$ws.Range("Q245").Value = "`$/caja 18 kilos"
$ws.Range("R245").Value = "Provincia de Limarí"
$ws.Range("S245").Value = 556
$ws.Range("T245").Value = 18

# New row 246
$ws.Range("A246").Value = 5
$ws.Range("B246").Value = "Macroferia Regional de Talca"
$ws.Range("C246").Value = "Maule"
$ws.Range("D246").Value = 44714
$ws.Range("E246").Value = 7
$ws.Range("F246").Value = "Fruta"
$ws.Range("G246").Value = 100102
$ws.Range("H246").Value = "Cítricos"
$ws.Range("I246").Value = 100102004
$ws.Range("J246").Value = "Mandarina"
$ws.Range("K246").Value = "Clemenuless"
$ws.Range("L246").Value = "Segunda"
$ws.Range("M246").Value = 200
$ws.Range("N246").Value = 5000
$ws.Range("O246").Value = 5000
$ws.Range("P246").Value = 5000
$ws.Range("Q246").Value = "`$/caja 10 kilos"
$ws.Range("R246").Value = "Provincia de Limarí"
$ws.Range("S246").Value = 500
$ws.Range("T246").Value = 10
